$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Add the new "GuestUser" sheet right after "AddProduct"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "GuestUser"

# 2) Pre-apply the cell formats (border-only for data rows, bold+border for
#    the header row) by copying them from the AddProduct sheet, so the new
#    styles line up with (and re-use) the existing style table entries
#    instead of minting new ones.
$ws1.Range("A1").Copy()
$ws2.Range("A1:S1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:S2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Header row (row 1) - fill in the exact order the field names were
#    originally entered (ConfirmPasswd's string is registered before
#    Passwd's, even though Passwd's column comes first).
$ws2.Cells.Item(1,1).Value  = "Email"
$ws2.Cells.Item(1,2).Value  = "UserName"
$ws2.Cells.Item(1,3).Value  = "DisplayedName"
$ws2.Cells.Item(1,5).Value  = "ConfirmPasswd"
$ws2.Cells.Item(1,4).Value  = "Passwd"
$ws2.Cells.Item(1,6).Value  = "CompanyName"
$ws2.Cells.Item(1,7).Value  = "Title"
$ws2.Cells.Item(1,8).Value  = "FirstName"
$ws2.Cells.Item(1,9).Value  = "MiddleName"
$ws2.Cells.Item(1,10).Value = "LastName"
$ws2.Cells.Item(1,11).Value = "Address1"
$ws2.Cells.Item(1,12).Value = "Address2"
$ws2.Cells.Item(1,13).Value = "PostalCode"
$ws2.Cells.Item(1,14).Value = "City"
$ws2.Cells.Item(1,15).Value = "Country"
$ws2.Cells.Item(1,16).Value = "State/Province"
$ws2.Cells.Item(1,17).Value = "Phone"
$ws2.Cells.Item(1,18).Value = "MobilePhone"
$ws2.Cells.Item(1,19).Value = "Fax"

# 4) Data row (row 2) - sample guest-user record. Values that look like
#    numbers but must stay text (so leading zeros / full digit strings
#    survive) are entered with a leading apostrophe, matching how the
#    workbook stores them (shared string + quotePrefix style).
$ws2.Cells.Item(2,1).Value  = "thohip@gmail.com"
$ws2.Cells.Item(2,2).Value  = "ThoHip"
$ws2.Cells.Item(2,3).Value  = "ThoHip"
$ws2.Cells.Item(2,4).Value  = 12345678
$ws2.Cells.Item(2,5).Value  = "'12345678"
$ws2.Cells.Item(2,6).Value  = "HunterPro"
$ws2.Cells.Item(2,8).Value  = "Jack"
$ws2.Cells.Item(2,10).Value = "Phan"
$ws2.Cells.Item(2,11).Value = "HCM123"
$ws2.Cells.Item(2,12).Value = "HCM123"
$ws2.Cells.Item(2,13).Value = "'00008"
$ws2.Cells.Item(2,14).Value = "HCM"
$ws2.Cells.Item(2,15).Value = "VietName"
$ws2.Cells.Item(2,16).Value = "HCM"
$ws2.Cells.Item(2,17).Value = "'0979155626"
$ws2.Cells.Item(2,18).Value = "'0979155626"
$ws2.Cells.Item(2,19).Value = 83123456
$ws2.Cells.Item(2,7).Value  = "Mrs"

# 5) Column widths (best-fit widths as authored; the engine quantizes
#    ColumnWidth to 1/6-character steps internally, so these inputs are
#    chosen to land on the nearest achievable stored width to the original
#    bestFit values)
$ws2.Columns.Item(1).ColumnWidth = 17.25
$ws2.Columns.Item(2).ColumnWidth = 9.5834
$ws2.Columns.Item(3).ColumnWidth = 14.25
$ws2.Columns.Item(4).ColumnWidth = 8.0834
$ws2.Columns.Item(5).ColumnWidth = 13.9167
$ws2.Columns.Item(6).ColumnWidth = 13.75
$ws2.Columns.Item(7).ColumnWidth = 4.0834
$ws2.Columns.Item(8).ColumnWidth = 9.25
$ws2.Columns.Item(9).ColumnWidth = 11.75
$ws2.Columns.Item(10).ColumnWidth = 8.75
$ws2.Columns.Item(13).ColumnWidth = 10.0834
$ws2.Columns.Item(14).ColumnWidth = 3.5834
$ws2.Columns.Item(15).ColumnWidth = 7.0834
$ws2.Columns.Item(16).ColumnWidth = 13.4167
$ws2.Columns.Item(17).ColumnWidth = 5.75
$ws2.Columns.Item(18).ColumnWidth = 12.4167
$ws2.Columns.Item(19).ColumnWidth = 8.0834

# 6) View tweaks
$ws2.Application.ActiveWindow.ScrollColumn = 2
$ws2.Range("I15").Select()

# 7) Sheet1 tweaks: selection moved, and the two "quantity-like" text cells
#    (C2/C3) switch to an explicit Text number format (keeps them stored as
#    shared strings with quotePrefix, same values, new style entry).
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C8").Select()
